# Update "paises" (countries) COVID-19 dashboard data and the
# "provincias Spain" (Spain provinces) refresh timestamp.
#
# The underlying source data was refreshed: case totals grew for a
# number of countries, and three pairs of neighbouring rows swapped
# rank order (the table is kept sorted by "Casos totales" descending),
# so the country name together with its row of figures moves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh timestamp shown in row 1.
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 12:21"

# Row number -> Pais, Casos totales, Nuevos casos, Casos activos,
#               Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @{
    4   = @("Estados Unidos",        3695581, 556,  1680418, 1874038, 0, 7,   141125)
    6   = @("India",                 1008480, 2843, 637646,  345170,  0, 55,  25664)
    14  = @("Iran",                  269440,  2379, 232873,  22776,   0, 183, 13791)
    29  = @("Indonesia",             83130,   1462, 41834,   37339,   0, 84,  3957)
    34  = @("Oman",                  64193,   1619, 41450,   22445,   0, 8,   298)
    35  = @("Belgica",               63238,   199,  17253,   36190,   0, 3,   9795)
    36  = @("Filipinas",             63001,   1841, 21748,   39593,   0, 17,  1660)
    48  = @("Rumania",               35802,   799,  22312,   11502,   0, 17,  1988)
    49  = @("Afganistan",            35229,   159,  23151,   10931,   0, 34,  1147)
    50  = @("Barein",                35084,   0,    30809,   4154,    0, 0,   121)
    65  = @("Marruecos",             16638,   93,   14175,   2200,    0, 0,   263)
    74  = @("Australia",             11235,   425,  8117,    3002,    0, 3,   116)
    86  = @("Estado de Palestina",   7764,    352,  1492,    6221,    0, 0,   51)
    87  = @("Bosnia y Herzegovina",  7681,    0,    3534,    3907,    0, 0,   240)
    102 = @("Albania",               3906,    55,   2214,    1585,    0, 3,   107)
    141 = @("Uganda",                1056,    5,    1023,    33,      0, 0,   0)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("A$r").Value = $vals[0]
    $ws.Range("B$r").Value = $vals[1]
    $ws.Range("C$r").Value = $vals[2]
    $ws.Range("D$r").Value = $vals[3]
    $ws.Range("E$r").Value = $vals[4]
    $ws.Range("F$r").Value = $vals[5]
    $ws.Range("G$r").Value = $vals[6]
    $ws.Range("H$r").Value = $vals[7]
}
